$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$failures = @()
$result = $d.Content.Find.Execute("17+75=92", $false, $false, $false, $false, $false, $true, 1, $false, "67-27=40", 1)
if (-not $result) { $failures += "17+75=92" }
$result = $d.Content.Find.Execute("70+26=96", $false, $false, $false, $false, $false, $true, 1, $false, "44-6=38", 1)
if (-not $result) { $failures += "70+26=96" }
$result = $d.Content.Find.Execute("19+40=59", $false, $false, $false, $false, $false, $true, 1, $false, "41-23=18", 1)
if (-not $result) { $failures += "19+40=59" }
$result = $d.Content.Find.Execute("3+26=29", $false, $false, $false, $false, $false, $true, 1, $false, "22+76=98", 1)
if (-not $result) { $failures += "3+26=29" }
$result = $d.Content.Find.Execute("98-62=36", $false, $false, $false, $false, $false, $true, 1, $false, "88-75=13", 1)
if (-not $result) { $failures += "98-62=36" }
$result = $d.Content.Find.Execute("77-6=71", $false, $false, $false, $false, $false, $true, 1, $false, "64+27=91", 1)
if (-not $result) { $failures += "77-6=71" }
$result = $d.Content.Find.Execute("37+11=48", $false, $false, $false, $false, $false, $true, 1, $false, "83-11=72", 1)
if (-not $result) { $failures += "37+11=48" }
$result = $d.Content.Find.Execute("19+46=65", $false, $false, $false, $false, $false, $true, 1, $false, "34+38=72", 1)
if (-not $result) { $failures += "19+46=65" }
$result = $d.Content.Find.Execute("73+23=96", $false, $false, $false, $false, $false, $true, 1, $false, "68-38=30", 1)
if (-not $result) { $failures += "73+23=96" }
$result = $d.Content.Find.Execute("37+27=64", $false, $false, $false, $false, $false, $true, 1, $false, "55-54=1", 1)
if (-not $result) { $failures += "37+27=64" }
$result = $d.Content.Find.Execute("64+10=74", $false, $false, $false, $false, $false, $true, 1, $false, "53+37=90", 1)
if (-not $result) { $failures += "64+10=74" }
$result = $d.Content.Find.Execute("89-10=79", $false, $false, $false, $false, $false, $true, 1, $false, "34-19=15", 1)
if (-not $result) { $failures += "89-10=79" }
$result = $d.Content.Find.Execute("57+16=73", $false, $false, $false, $false, $false, $true, 1, $false, "10+73=83", 1)
if (-not $result) { $failures += "57+16=73" }
$result = $d.Content.Find.Execute("45+21=66", $false, $false, $false, $false, $false, $true, 1, $false, "35+8=43", 1)
if (-not $result) { $failures += "45+21=66" }
$result = $d.Content.Find.Execute("94-69=25", $false, $false, $false, $false, $false, $true, 1, $false, "87-15=72", 1)
if (-not $result) { $failures += "94-69=25" }
$result = $d.Content.Find.Execute("73-35=38", $false, $false, $false, $false, $false, $true, 1, $false, "43+26=69", 1)
if (-not $result) { $failures += "73-35=38" }
$result = $d.Content.Find.Execute("79-75=4", $false, $false, $false, $false, $false, $true, 1, $false, "86-5=81", 1)
if (-not $result) { $failures += "79-75=4" }
$result = $d.Content.Find.Execute("62+37=99", $false, $false, $false, $false, $false, $true, 1, $false, "29+37=66", 1)
if (-not $result) { $failures += "62+37=99" }
$result = $d.Content.Find.Execute("66-35=31", $false, $false, $false, $false, $false, $true, 1, $false, "47+35=82", 1)
if (-not $result) { $failures += "66-35=31" }
$result = $d.Content.Find.Execute("89-80=9", $false, $false, $false, $false, $false, $true, 1, $false, "86+7=93", 1)
if (-not $result) { $failures += "89-80=9" }
$result = $d.Content.Find.Execute("88-73=15", $false, $false, $false, $false, $false, $true, 1, $false, "51-41=10", 1)
if (-not $result) { $failures += "88-73=15" }
$result = $d.Content.Find.Execute("43+33=76", $false, $false, $false, $false, $false, $true, 1, $false, "1+92=93", 1)
if (-not $result) { $failures += "43+33=76" }
$result = $d.Content.Find.Execute("82-59=23", $false, $false, $false, $false, $false, $true, 1, $false, "47+32=79", 1)
if (-not $result) { $failures += "82-59=23" }
$result = $d.Content.Find.Execute("75-30=45", $false, $false, $false, $false, $false, $true, 1, $false, "25+64=89", 1)
if (-not $result) { $failures += "75-30=45" }
$result = $d.Content.Find.Execute("51-5=46", $false, $false, $false, $false, $false, $true, 1, $false, "29+66=95", 1)
if (-not $result) { $failures += "51-5=46" }
$result = $d.Content.Find.Execute("18+76=94", $false, $false, $false, $false, $false, $true, 1, $false, "85-26=59", 1)
if (-not $result) { $failures += "18+76=94" }
$result = $d.Content.Find.Execute("60-37=23", $false, $false, $false, $false, $false, $true, 1, $false, "63-34=29", 1)
if (-not $result) { $failures += "60-37=23" }
$result = $d.Content.Find.Execute("78-6=72", $false, $false, $false, $false, $false, $true, 1, $false, "63+9=72", 1)
if (-not $result) { $failures += "78-6=72" }
$result = $d.Content.Find.Execute("81+11=92", $false, $false, $false, $false, $false, $true, 1, $false, "55+14=69", 1)
if (-not $result) { $failures += "81+11=92" }
$result = $d.Content.Find.Execute("51-47=4", $false, $false, $false, $false, $false, $true, 1, $false, "21+24=45", 1)
if (-not $result) { $failures += "51-47=4" }
$result = $d.Content.Find.Execute("66-32=34", $false, $false, $false, $false, $false, $true, 1, $false, "4+69=73", 1)
if (-not $result) { $failures += "66-32=34" }
$result = $d.Content.Find.Execute("34+14=48", $false, $false, $false, $false, $false, $true, 1, $false, "54+14=68", 1)
if (-not $result) { $failures += "34+14=48" }
$result = $d.Content.Find.Execute("43+15=58", $false, $false, $false, $false, $false, $true, 1, $false, "74-13=61", 1)
if (-not $result) { $failures += "43+15=58" }
$result = $d.Content.Find.Execute("44-10=34", $false, $false, $false, $false, $false, $true, 1, $false, "80-7=73", 1)
if (-not $result) { $failures += "44-10=34" }
$result = $d.Content.Find.Execute("24-2=22", $false, $false, $false, $false, $false, $true, 1, $false, "54-23=31", 1)
if (-not $result) { $failures += "24-2=22" }
$result = $d.Content.Find.Execute("21+40=61", $false, $false, $false, $false, $false, $true, 1, $false, "11-7=4", 1)
if (-not $result) { $failures += "21+40=61" }
$result = $d.Content.Find.Execute("66+9=75", $false, $false, $false, $false, $false, $true, 1, $false, "52-28=24", 1)
if (-not $result) { $failures += "66+9=75" }
$result = $d.Content.Find.Execute("34+48=82", $false, $false, $false, $false, $false, $true, 1, $false, "15-0=15", 1)
if (-not $result) { $failures += "34+48=82" }
$result = $d.Content.Find.Execute("42-37=5", $false, $false, $false, $false, $false, $true, 1, $false, "58-25=33", 1)
if (-not $result) { $failures += "42-37=5" }
$result = $d.Content.Find.Execute("74-71=3", $false, $false, $false, $false, $false, $true, 1, $false, "73-65=8", 1)
if (-not $result) { $failures += "74-71=3" }
$result = $d.Content.Find.Execute("27-24=3", $false, $false, $false, $false, $false, $true, 1, $false, "39-27=12", 1)
if (-not $result) { $failures += "27-24=3" }
$result = $d.Content.Find.Execute("40+58=98", $false, $false, $false, $false, $false, $true, 1, $false, "59-47=12", 1)
if (-not $result) { $failures += "40+58=98" }
$result = $d.Content.Find.Execute("81-10=71", $false, $false, $false, $false, $false, $true, 1, $false, "45+27=72", 1)
if (-not $result) { $failures += "81-10=71" }
$result = $d.Content.Find.Execute("65-7=58", $false, $false, $false, $false, $false, $true, 1, $false, "35+25=60", 1)
if (-not $result) { $failures += "65-7=58" }
$result = $d.Content.Find.Execute("59-49=10", $false, $false, $false, $false, $false, $true, 1, $false, "75-43=32", 1)
if (-not $result) { $failures += "59-49=10" }
$result = $d.Content.Find.Execute("66+5=71", $false, $false, $false, $false, $false, $true, 1, $false, "31+66=97", 1)
if (-not $result) { $failures += "66+5=71" }
$result = $d.Content.Find.Execute("76+19=95", $false, $false, $false, $false, $false, $true, 1, $false, "4+64=68", 1)
if (-not $result) { $failures += "76+19=95" }
$result = $d.Content.Find.Execute("92-2=90", $false, $false, $false, $false, $false, $true, 1, $false, "19+7=26", 1)
if (-not $result) { $failures += "92-2=90" }
$result = $d.Content.Find.Execute("52+24=76", $false, $false, $false, $false, $false, $true, 1, $false, "38+35=73", 1)
if (-not $result) { $failures += "52+24=76" }
$result = $d.Content.Find.Execute("3+40=43", $false, $false, $false, $false, $false, $true, 1, $false, "69+17=86", 1)
if (-not $result) { $failures += "3+40=43" }
$result = $d.Content.Find.Execute("49-42=7", $false, $false, $false, $false, $false, $true, 1, $false, "42+23=65", 1)
if (-not $result) { $failures += "49-42=7" }
$result = $d.Content.Find.Execute("84-30=54", $false, $false, $false, $false, $false, $true, 1, $false, "85+8=93", 1)
if (-not $result) { $failures += "84-30=54" }
$result = $d.Content.Find.Execute("57+28=85", $false, $false, $false, $false, $false, $true, 1, $false, "95+2=97", 1)
if (-not $result) { $failures += "57+28=85" }
$result = $d.Content.Find.Execute("27+26=53", $false, $false, $false, $false, $false, $true, 1, $false, "91-45=46", 1)
if (-not $result) { $failures += "27+26=53" }
$result = $d.Content.Find.Execute("48+22=70", $false, $false, $false, $false, $false, $true, 1, $false, "64-44=20", 1)
if (-not $result) { $failures += "48+22=70" }
$result = $d.Content.Find.Execute("48+4=52", $false, $false, $false, $false, $false, $true, 1, $false, "87-62=25", 1)
if (-not $result) { $failures += "48+4=52" }
$result = $d.Content.Find.Execute("34+9=43", $false, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 1)
if (-not $result) { $failures += "34+9=43" }
$result = $d.Content.Find.Execute("59-19=40", $false, $false, $false, $false, $false, $true, 1, $false, "46-0=46", 1)
if (-not $result) { $failures += "59-19=40" }
$result = $d.Content.Find.Execute("60+39=99", $false, $false, $false, $false, $false, $true, 1, $false, "5+16=21", 1)
if (-not $result) { $failures += "60+39=99" }
$result = $d.Content.Find.Execute("6+70=76", $false, $false, $false, $false, $false, $true, 1, $false, "90-47=43", 1)
if (-not $result) { $failures += "6+70=76" }
$result = $d.Content.Find.Execute("76+15=91", $false, $false, $false, $false, $false, $true, 1, $false, "79-42=37", 1)
if (-not $result) { $failures += "76+15=91" }
$result = $d.Content.Find.Execute("70-8=62", $false, $false, $false, $false, $false, $true, 1, $false, "82-16=66", 1)
if (-not $result) { $failures += "70-8=62" }
$result = $d.Content.Find.Execute("61+13=74", $false, $false, $false, $false, $false, $true, 1, $false, "60-32=28", 1)
if (-not $result) { $failures += "61+13=74" }
$result = $d.Content.Find.Execute("15+49=64", $false, $false, $false, $false, $false, $true, 1, $false, "24+11=35", 1)
if (-not $result) { $failures += "15+49=64" }
$result = $d.Content.Find.Execute("69-56=13", $false, $false, $false, $false, $false, $true, 1, $false, "38+6=44", 1)
if (-not $result) { $failures += "69-56=13" }
$result = $d.Content.Find.Execute("76-24=52", $false, $false, $false, $false, $false, $true, 1, $false, "53-22=31", 1)
if (-not $result) { $failures += "76-24=52" }
$result = $d.Content.Find.Execute("43+10=53", $false, $false, $false, $false, $false, $true, 1, $false, "17+74=91", 1)
if (-not $result) { $failures += "43+10=53" }
$result = $d.Content.Find.Execute("19+56=75", $false, $false, $false, $false, $false, $true, 1, $false, "68-23=45", 1)
if (-not $result) { $failures += "19+56=75" }
$result = $d.Content.Find.Execute("5+67=72", $false, $false, $false, $false, $false, $true, 1, $false, "83-38=45", 1)
if (-not $result) { $failures += "5+67=72" }
$result = $d.Content.Find.Execute("38-26=12", $false, $false, $false, $false, $false, $true, 1, $false, "50+18=68", 1)
if (-not $result) { $failures += "38-26=12" }
$result = $d.Content.Find.Execute("98-82=16", $false, $false, $false, $false, $false, $true, 1, $false, "41+35=76", 1)
if (-not $result) { $failures += "98-82=16" }
$result = $d.Content.Find.Execute("2+70=72", $false, $false, $false, $false, $false, $true, 1, $false, "26+14=40", 1)
if (-not $result) { $failures += "2+70=72" }
$result = $d.Content.Find.Execute("71-12=59", $false, $false, $false, $false, $false, $true, 1, $false, "83-31=52", 1)
if (-not $result) { $failures += "71-12=59" }
$result = $d.Content.Find.Execute("2+56=58", $false, $false, $false, $false, $false, $true, 1, $false, "65-9=56", 1)
if (-not $result) { $failures += "2+56=58" }
$result = $d.Content.Find.Execute("94-86=8", $false, $false, $false, $false, $false, $true, 1, $false, "30+18=48", 1)
if (-not $result) { $failures += "94-86=8" }
$result = $d.Content.Find.Execute("52+43=95", $false, $false, $false, $false, $false, $true, 1, $false, "59-7=52", 1)
if (-not $result) { $failures += "52+43=95" }
$result = $d.Content.Find.Execute("10+67=77", $false, $false, $false, $false, $false, $true, 1, $false, "12+5=17", 1)
if (-not $result) { $failures += "10+67=77" }
$result = $d.Content.Find.Execute("16+10=26", $false, $false, $false, $false, $false, $true, 1, $false, "25+25=50", 1)
if (-not $result) { $failures += "16+10=26" }
$result = $d.Content.Find.Execute("49-3=46", $false, $false, $false, $false, $false, $true, 1, $false, "86-71=15", 1)
if (-not $result) { $failures += "49-3=46" }
$result = $d.Content.Find.Execute("36-28=8", $false, $false, $false, $false, $false, $true, 1, $false, "8+19=27", 1)
if (-not $result) { $failures += "36-28=8" }
$result = $d.Content.Find.Execute("68+14=82", $false, $false, $false, $false, $false, $true, 1, $false, "90-40=50", 1)
if (-not $result) { $failures += "68+14=82" }
$result = $d.Content.Find.Execute("89-19=70", $false, $false, $false, $false, $false, $true, 1, $false, "9+10=19", 1)
if (-not $result) { $failures += "89-19=70" }
$result = $d.Content.Find.Execute("72-69=3", $false, $false, $false, $false, $false, $true, 1, $false, "71+12=83", 1)
if (-not $result) { $failures += "72-69=3" }
$result = $d.Content.Find.Execute("78-69=9", $false, $false, $false, $false, $false, $true, 1, $false, "54-19=35", 1)
if (-not $result) { $failures += "78-69=9" }
$result = $d.Content.Find.Execute("93-73=20", $false, $false, $false, $false, $false, $true, 1, $false, "91-65=26", 1)
if (-not $result) { $failures += "93-73=20" }
$result = $d.Content.Find.Execute("92-73=19", $false, $false, $false, $false, $false, $true, 1, $false, "6+78=84", 1)
if (-not $result) { $failures += "92-73=19" }
$result = $d.Content.Find.Execute("38-12=26", $false, $false, $false, $false, $false, $true, 1, $false, "79-36=43", 1)
if (-not $result) { $failures += "38-12=26" }
$result = $d.Content.Find.Execute("48-19=29", $false, $false, $false, $false, $false, $true, 1, $false, "82+6=88", 1)
if (-not $result) { $failures += "48-19=29" }
$result = $d.Content.Find.Execute("62+6=68", $false, $false, $false, $false, $false, $true, 1, $false, "43-35=8", 1)
if (-not $result) { $failures += "62+6=68" }
$result = $d.Content.Find.Execute("31+57=88", $false, $false, $false, $false, $false, $true, 1, $false, "36+60=96", 1)
if (-not $result) { $failures += "31+57=88" }
$result = $d.Content.Find.Execute("73+4=77", $false, $false, $false, $false, $false, $true, 1, $false, "16+56=72", 1)
if (-not $result) { $failures += "73+4=77" }
$result = $d.Content.Find.Execute("48+4=52", $false, $false, $false, $false, $false, $true, 1, $false, "16+0=16", 1)
if (-not $result) { $failures += "48+4=52" }
$result = $d.Content.Find.Execute("64-4=60", $false, $false, $false, $false, $false, $true, 1, $false, "16+3=19", 1)
if (-not $result) { $failures += "64-4=60" }
$result = $d.Content.Find.Execute("82-32=50", $false, $false, $false, $false, $false, $true, 1, $false, "16+44=60", 1)
if (-not $result) { $failures += "82-32=50" }
$result = $d.Content.Find.Execute("1+63=64", $false, $false, $false, $false, $false, $true, 1, $false, "92-66=26", 1)
if (-not $result) { $failures += "1+63=64" }
$result = $d.Content.Find.Execute("63-44=19", $false, $false, $false, $false, $false, $true, 1, $false, "48-38=10", 1)
if (-not $result) { $failures += "63-44=19" }
$result = $d.Content.Find.Execute("78-3=75", $false, $false, $false, $false, $false, $true, 1, $false, "92-86=6", 1)
if (-not $result) { $failures += "78-3=75" }
$result = $d.Content.Find.Execute("0+47=47", $false, $false, $false, $false, $false, $true, 1, $false, "15+56=71", 1)
if (-not $result) { $failures += "0+47=47" }
$result = $d.Content.Find.Execute("6+80=86", $false, $false, $false, $false, $false, $true, 1, $false, "32+1=33", 1)
if (-not $result) { $failures += "6+80=86" }
$result = $d.Content.Find.Execute("50+29=79", $false, $false, $false, $false, $false, $true, 1, $false, "33+26=59", 1)
if (-not $result) { $failures += "50+29=79" }
Write-Output "Failures: $($failures.Count)"
if ($failures.Count -gt 0) { Write-Output ($failures -join ", ") }
